# Auto-generated Excel COM-interop script to update the cryptos table
# per the commit "Updated cryptos list on Mon Jun 19 18:29:50 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "'26.800.03"
$ws.Cells.Item(2,5).Value = "  +0.74%  "

# Row 3
$ws.Cells.Item(3,4).Value = "'1.733.91"
$ws.Cells.Item(3,5).Value = "  -0.33%  "

# Row 4
$ws.Cells.Item(4,4).Value = "'0.9950"
$ws.Cells.Item(4,5).Value = "  -0.53%  "

# Row 5
$ws.Cells.Item(5,4).Value = "'242.46"
$ws.Cells.Item(5,5).Value = "  -1.53%  "

# Row 6
$ws.Cells.Item(6,4).Value = "'0.9959"
$ws.Cells.Item(6,5).Value = "  -0.46%  "

# Row 7
$ws.Cells.Item(7,4).Value = "'0.4934"
$ws.Cells.Item(7,5).Value = "  +0.55%  "

# Row 8
$ws.Cells.Item(8,4).Value = "'0.2614"
$ws.Cells.Item(8,5).Value = "  -2.27%  "

# Row 9
$ws.Cells.Item(9,4).Value = "'0.06232"
$ws.Cells.Item(9,5).Value = "  -0.78%  "

# Row 10
$ws.Cells.Item(10,4).Value = "'1.724.95"
$ws.Cells.Item(10,5).Value = "  -0.94%  "

# Row 11
$ws.Cells.Item(11,4).Value = "'15.81"
$ws.Cells.Item(11,5).Value = "  +0.64%  "

# Row 12
$ws.Cells.Item(12,4).Value = "'0.06981"
$ws.Cells.Item(12,5).Value = "  -0.95%  "

# Row 13
$ws.Cells.Item(13,4).Value = "'0.6168"
$ws.Cells.Item(13,5).Value = "  +0.54%  "

# Row 14
$ws.Cells.Item(14,4).Value = "'4.513"
$ws.Cells.Item(14,5).Value = "  -1.59%  "

# Row 15
$ws.Cells.Item(15,4).Value = "'77.25"
$ws.Cells.Item(15,5).Value = "  -0.98%  "

# Row 16
$ws.Cells.Item(16,4).Value = "'0.9950"
$ws.Cells.Item(16,5).Value = "  -0.54%  "

# Row 17
$ws.Cells.Item(17,4).Value = "'26.548.40"
$ws.Cells.Item(17,5).Value = "  -0.24%  "

# Row 18
$ws.Cells.Item(18,4).Value = "'0.9946"
$ws.Cells.Item(18,5).Value = "  -0.63%  "

# Row 19
$ws.Cells.Item(19,4).Value = "'0.000007184"
$ws.Cells.Item(19,5).Value = "  -1.42%  "

# Row 20
$ws.Cells.Item(20,4).Value = "'11.44"
$ws.Cells.Item(20,5).Value = "  -1.10%  "

# Row 21
$ws.Cells.Item(21,4).Value = "'1.949.23"
$ws.Cells.Item(21,5).Value = "  -1.24%  "

# Row 22
$ws.Cells.Item(22,4).Value = "'4.459"
$ws.Cells.Item(22,5).Value = "  -2.34%  "

# Row 23
$ws.Cells.Item(23,4).Value = "'8.548"
$ws.Cells.Item(23,5).Value = "  -1.89%  "

# Row 24
$ws.Cells.Item(24,4).Value = "'5.153"
$ws.Cells.Item(24,5).Value = "  -2.30%  "

# Row 25
$ws.Cells.Item(25,4).Value = "'138.70"
$ws.Cells.Item(25,5).Value = "  -0.17%  "

# Row 26
$ws.Cells.Item(26,4).Value = "'15.36"
$ws.Cells.Item(26,5).Value = "  -0.39%  "

# Row 27
$ws.Cells.Item(27,4).Value = "'1.411"
$ws.Cells.Item(27,5).Value = "  -0.69%  "

# Row 28
$ws.Cells.Item(28,4).Value = "'1.762"
$ws.Cells.Item(28,5).Value = "  +0.01%  "

# Row 29
$ws.Cells.Item(29,4).Value = "'106.63"
$ws.Cells.Item(29,5).Value = "  -0.72%  "

# Row 30
$ws.Cells.Item(30,4).Value = "'3.947"
$ws.Cells.Item(30,5).Value = "  -1.96%  "

# Row 31
$ws.Cells.Item(31,4).Value = "'0.07981"
$ws.Cells.Item(31,5).Value = "  -0.92%  "

# Row 32
$ws.Cells.Item(32,4).Value = "'3.668"
$ws.Cells.Item(32,5).Value = "  -1.53%  "

# Row 33
$ws.Cells.Item(33,4).Value = "'0.04528"
$ws.Cells.Item(33,5).Value = "  -1.70%  "

# Row 34
$ws.Cells.Item(34,2).Value = "Frax"
$ws.Cells.Item(34,3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(34,4).Value = "'0.9946"
$ws.Cells.Item(34,5).Value = "  -0.52%  "

# Row 35
$ws.Cells.Item(35,2).Value = "HuobiToken"
$ws.Cells.Item(35,3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35,4).Value = "'2.614"
$ws.Cells.Item(35,5).Value = "  -0.05%  "

# Row 36
$ws.Cells.Item(36,2).Value = "ARBITRUM"
$ws.Cells.Item(36,3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36,4).Value = "'1.007"
$ws.Cells.Item(36,5).Value = "  -0.61%  "

# Row 37
$ws.Cells.Item(37,2).Value = "ImmutableX"
$ws.Cells.Item(37,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(37,4).Value = "'0.6290"
$ws.Cells.Item(37,5).Value = "  -1.71%  "

# Row 38
$ws.Cells.Item(38,2).Value = "TrustWalletToken"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(38,4).Value = "'0.9524"
$ws.Cells.Item(38,5).Value = "  +5.22%  "

# Row 39
$ws.Cells.Item(39,2).Value = "RenderToken"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39,4).Value = "'2.019"
$ws.Cells.Item(39,5).Value = "  -2.45%  "

# Row 40
$ws.Cells.Item(40,2).Value = "MXToken"
$ws.Cells.Item(40,3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(40,4).Value = "'2.424"
$ws.Cells.Item(40,5).Value = "  -0.25%  "

# Row 41
$ws.Cells.Item(41,2).Value = "PaxDollar"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Cells.Item(41,4).Value = "'0.9949"
$ws.Cells.Item(41,5).Value = "  -0.81%  "

# Row 42
$ws.Cells.Item(42,2).Value = "VeChain"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(42,4).Value = "'0.01507"
$ws.Cells.Item(42,5).Value = "  +0.26%  "

# Row 43
$ws.Cells.Item(43,2).Value = "Quant"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(43,4).Value = "'99.72"
$ws.Cells.Item(43,5).Value = "  -2.39%  "

# Row 44
$ws.Cells.Item(44,2).Value = "FraxShare"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(44,4).Value = "'5.497"
$ws.Cells.Item(44,5).Value = "  +1.21%  "

# Row 45
$ws.Cells.Item(45,2).Value = "TheSandbox"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(45,4).Value = "'0.3880"
$ws.Cells.Item(45,5).Value = "  -0.88%  "

# Row 46
$ws.Cells.Item(46,2).Value = "Aptos"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(46,4).Value = "'6.962"
$ws.Cells.Item(46,5).Value = "  +1.37%  "

# Row 47
$ws.Cells.Item(47,2).Value = "Algorand"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(47,4).Value = "'0.1165"
$ws.Cells.Item(47,5).Value = "  -1.59%  "

# Row 48
$ws.Cells.Item(48,2).Value = "Cronos"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(48,4).Value = "'0.05398"
$ws.Cells.Item(48,5).Value = "  +0.03%  "

# Row 49
$ws.Cells.Item(49,2).Value = "Elrond"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(49,4).Value = "'30.58"
$ws.Cells.Item(49,5).Value = "  -0.09%  "

# Row 50
$ws.Cells.Item(50,2).Value = "EnergySwap"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50,4).Value = "'7.800"
$ws.Cells.Item(50,5).Value = "  -0.03%  "

# Row 51
$ws.Cells.Item(51,2).Value = "Aave"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(51,4).Value = "'51.86"
$ws.Cells.Item(51,5).Value = "  +0.14%  "

